$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.857.34"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.08"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.60"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.63"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.27"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.23"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.588.71"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.333.66"
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.62"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.068.71"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0970"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.15"
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.51"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.33"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.69"
$ws.Range("E29").Value = "  +5.34%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.20"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.94"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.75"
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0800"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("E36").Value = "  +3.89%  "
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.49"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.61"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.746.97"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.91"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.30"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.95"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.15"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.59"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.01"
$ws.Range("E51").Value = "  -1.50%  "
